# Edit 1: "... change in a corresponding dependent variable. However, how ..."
#      -> "... change in an independent variable. But how ..."
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "change in a corresponding dependent variable. However, how certain",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "change in an independent variable. But how certain", 2)

# Edit 2: "... sample size and therefore the sample variation ..."
#      -> "... sample size as by increasing sample size we should be increasing the sample variation ..."
$d.Content.Find.Execute(
    "sample size and therefore the sample variation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sample size as by increasing sample size we should be increasing the sample variation", 2)
